$d = $word.ActiveDocument

# --- Step 1: append a run containing a single space right after "Tugas 02",
# as a run that stays separate from the "Tugas 02" run. A temporary bookmark
# dropped exactly at the split point stops the engine from coalescing the
# two otherwise identically-formatted runs when it is written back out; the
# bookmark is removed again immediately afterwards (without re-merging the
# now-separate runs). ---
$p1 = $d.Paragraphs(1)
$splitPos = $p1.Range.End - 1
$d.Bookmarks.Add("TempSplit", $d.Range($splitPos, $splitPos))
$d.Range($splitPos, $splitPos).InsertAfter(" ")
$d.Bookmarks("TempSplit").Delete()

# --- Step 2: insert a new paragraph "TIM : KAMI Kembali" right after the
# "Tugas 02" paragraph. ---
$p1 = $d.Paragraphs(1)
$endOfP1 = $p1.Range.End - 1
$d.Range($endOfP1, $endOfP1).InsertParagraphAfter()

$p2 = $d.Paragraphs(2)
$p2.Range.InsertBefore("TIM : KAMI Kembali")

# --- Step 3: relocate the hidden "_GoBack" bookmark (Word's "last edit"
# marker) from wherever it currently sits to the end of the text just
# typed ("TIM : KAMI Kembali"), matching where a real edit session would
# leave it. A placeholder trailing character keeps the bookmark's range
# off the exact paragraph-end boundary (collapsed ranges sitting exactly
# on a paragraph boundary resolve incorrectly), and is removed right
# after, leaving the bookmark collapsed at the true end of the text. ---
$p2 = $d.Paragraphs(2)
$placeholderPos = $p2.Range.End - 1
$d.Range($placeholderPos, $placeholderPos).InsertAfter("#")
$d.Bookmarks.Add("_GoBack", $d.Range($placeholderPos, $placeholderPos))
$d.Range($placeholderPos, $placeholderPos + 1).Delete()
